# Refresh the cryptos table with the latest scraped Price/Volume(1h)
# figures, matching the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.402.30'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.01%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.329.01'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  -1.78%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.04'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -1.42%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.19'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -3.84%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.512'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -0.86%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.515'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -0.70%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.31'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -2.36%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0800'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -1.44%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '51.45'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -3.72%  '

# Row 13
$ws.Range("E13").Value = '  +0.61%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.84'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -2.18%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.683.56'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -2.06%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.82'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +1.30%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.312.55'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -2.63%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.808'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.251.07'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -0.24%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.89'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -0.95%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0907'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -1.38%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.12'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -2.95%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.72'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -0.83%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '238.43'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -1.21%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.99'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -3.16%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.54'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -2.96%  '

# Row 27
$ws.Range("E27").Value = '  +0.16%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.18'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -2.23%  '

# Row 29
$ws.Range("E29").Value = '  -5.56%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.90'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -5.14%  '

# Row 31
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.25'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -3.15%  '

# Row 32
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '165.42'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +2.38%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -0.17%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.09'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -3.51%  '

# Row 35
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.54'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -4.64%  '

# Row 36
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.42'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -4.72%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.97'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -7.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0711'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -4.20%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.90'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -6.37%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.83'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -5.72%  '

# Row 41
$ws.Range("E41").Value = '  -3.53%  '

# Row 42
$ws.Range("E42").Value = '  -2.20%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.44'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -9.10%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.984.25'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -1.62%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0286'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -1.39%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.73'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -5.13%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.96'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -5.80%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.91'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -6.57%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.16'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -4.97%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.86'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +3.22%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.550.68'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -0.56%  '
